$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Update the cached/serialized CDM add-in metadata stored in the A1 cell comment.
$comment = $ws.Range("A1").Comment
$newCommentBlob = "8RwAAB+LCAAAAAAAAAOlWVtvI0kV/istP8GD3W07M7mo0itfkmDhSxQ7ZLIvqN1diYu0u0xXdRK/LQK0o4WRkNgdBOyCQEIgJHZXYlmxO6D5L6txZnjiL3Dq0jfbmXEPo2jSdW516tSpc76qoHdup75xjUNGaLBfqlaskoEDl3okuNwvRfyiXH1YesdGB7cu9o+d0JliDsIGaAVs75aR/dKE89mead7c3FRu6hUaXpo1y6qaj3rdoTvBU6dMAsadwMWlRMt7s1bJRi1v2sPc8RzuKM39UmfYqbQwcdtA6zmBc4nDSjNiJMCMHQSccIKZ0Ayxw3Gr3fueWphdqzysVJG5Qk8lmxHxPSWXk1R0LQfT4hGZYrtmVXfK1m7Zqo+s7b1abW/rYWV7q/purJgIoq7D+BCH18SVhCF3pjOpbu1adWu7Vtt6iMy1QmArDYCNBr53gq8Jw14L+z4rFBFTb2DD5bDqYsG0kJnR1Ybe3oWj0JlNRoT7uJgbJ6dNYxxoX1IjNjqkIXYhfm/lUh/fDEId1tGsC9zRhIR83nbmhW2dMhwOZiJIxVRt1KYBb/g45Kcz2GvsQSoAw+ZhhJF5DzNVahPmwjcJIuzZF47Psko5Jjqj4RWbOS7uwzk2hY2bwKeOBwnHCePETSddYaDjkM7AIkzepL53CFa18BpGYrkTQIjFtE1Kr1Lv1jGR3FW5v7CnU4fH4it0NJzQm0Hgz4fRmLkhGWOv3Yyl1/KQOJBauxUxTqfgRUpCipah9HrmHP7BGVzmoDZ2ydTxj32II7PrYChHQI2I0wvCW9SPpgGL3VqiojNY1AjfJotMxmgA+xuIuNOgE8TyKtJrWXmFE3qTzLnKkHHIkBvMjXd8lbEs3AZavIOrHLkpYpWHxIcWkd2ODDWfGMMJxnxtVigOEtXwUDQduznvR9MxnLAxHLNrOStDZspHkKqQ7uCXbUEjKcufkWXtyR/wI2Gjg8C7Xy5mIpguM5ddBd4SCcGa/KbvBFdAPSN80m/Ea1nDQSoC98qv8hAc3pnvzCU5iVKWhjqB60ceVjWhE1zIFBW+qU29l41WSF045jZygvloPoPSzMgeh4/9EjTrPcZDgAMl26VRwMO5KB7I1KJv0mHROJATOP7GOhch/mEEKGR+GAVui3qbz+ap6JwGhG/uIY1CVRE3V5HRE8UxYm0syoys+xvru0XWxMJC4tMAT2lA3M2jDUEW3ntvsRAWn6qNNbA6XxvL+9DaVecTZ31jtRAgJPS6QtM0GKMukcmqj4eX0TfvOTJtfOFEPsA3Dl32Mqm9y2TUYFfLMlkSOg39uALaAhwzQMeuN624gB8EAqy4dCoIJoDSsyEys/ICBLn4ILjsOsFlBDAjqSvL9KT+ihY5Cp2AieUkqGKpFK8XQnGdUmjHVsVrEMlEUMWLAheZS3JohKczGjp+DwJDDnXaacgEaKTn8IkeQW/zsRsH2UxVE628Z7HjbxKTTUotQxx4XSaXiFJIrEXh8FQmpSGxyh4cS7/l+GQcqqoat/J1PNiwFB/G9VcsriBWjPcA7mLQfb+L5wKgpwNNlylbjRkqgUUhtYcnWzu1rd0dUJJjJFfcAp+oTyS+NI4oXAeDKegZgIwxIMc9Y87hKKjg6AZUQC+rgQ7j2i6n70GK5gh5PkCMSwKdY0UuYaTydg+A58SfZyTV6rrUBcHFLx4vvvry7sPni68/vvvwy8VP//Tff/1m8dUXi8efvHr/7y+efaDWp4TRyBn7WDo0au7sWPUtSK2EhEQ8TQmIvcjlknZ+LnFwMkb6OicHrYNO66jblCUkIcbqqouY4qY4p1E6HKpFyInkLprx5isRexSXJD3OcTNdyRY3t2ucl87y71NUsXj57M8vn/3tXm0dsBReVXd365ugL7gDb5Wt2j3oq5sr+1r4QblWywgvyaATVfOTOHU8u161dqs721Y1Kd9eksHrhJZZ2tLIuTSX9BSppRBRkgLZccyUeT/CjCdsdRIyAxXEu3/87NVnH+WkdHQ1JW8FnJPIRUxmxgNpun8yMoaD05PWgTE6GIo8SXkZOWX8NcJ69uQ8ZZNm9bTlUmodVZl7reLyjEchjWYrlSKlrpFcWzNWuUuaMjarrqW8NfJ61z5+vE5Br6SdAtLkCSNLQzmOImX4unY9+dGLr5+/eP6Tl//+1d3TzxZP/rj4549zdvRsyf0eUhcOSHaYZDJUMd01lijobChjemV9P9MlNFFch44pCTizaw+25FVIDxHoVoU5+Rt1ptC5pGUZNqAvUdB3HHZwy/VhtfvIzBPA0ZkDTZOmV8iEoOpyGt7/fPK7u99+cff081fv/3XxwV8WP3/68tnvX336B3WS7j76/O7Jp7pyLxd36Yu4mCosZ8iXDtcQJ8wQLdj45r1fGgHlBiAHI5JV5pv3fp0xJhyVGCO1DMgscSTvwopoVlnoGRlXEh9yeomK6uMt0ZbqiYRuTHRG3HSSd8vCFMBLQzK+1RmVI4YNCqDo27CSvHCqvKmeVlFt8njbqlVrmqu8EUsYOywT+iOfjh3fiBnyHWFJJKf1eoVUVs531B00G91URDkxCD0cijRUHyhGhqJNdFg8ilMtQwEu4Dc38sXTz4rYKiuxnClnpn5DuWh44kF9/YtDTgK1ojBUKCfQr/LDaAaYNn5su58vHyAzMLavIGcW2KbjTjvPh3GGC80tzxYEyZcVSrNUteow8VyjUGlfhCYdAi/3aAnh0A/vCj1dE9gSUxSegzCk4drqk3JisR4AYqgoZhrxREbuqQLPXrpXMSGueMmHusDpFdI29jEv9iptpto9QL9vqwt7X1S1wwa+p4NZ7AaRhCU1kH2aF4ny/77Mq2RrhCGAJfGOV/gpPb5/nsC1taA3ailSUVzkYHb91n1IQsYfiUqgvxTlPKGcK9T5SNyb1Iccn9v1B4oAAmbWuplzMz66XP2BhPpdMiUFb3dWfL7zRiCWs5mCZZ1imSJaSx/fAmjMWICiOP4BtA31HFLEmkpYqKWJvniCZORywos6tj12sIfHVtkd41p5y7N2yrsY18vVKvzvuLWaZT0QD5jaOFQOgm8KTmLGG5b+4dL+H2OnMenxHAAA=="
$comment.Text($newCommentBlob)

# 3. Cell text update: "Function Description" -> "Function Information"
$ws.Range("A11").Value = "Function Information"

# 4. Skewness value tweak (tiny precision correction)
$ws.Range("B20").Value = -0.3785388757796309

# 5. Number format for the observation values (B27:B36) changes from "0.000" to "###0.000"
$ws.Range("B27:B36").NumberFormat = "###0.000"
